# Removed implicit waits in the testcases and updated listeners
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update EXECUTE column (B) values
$ws.Range("B4").Value = "No"
$ws.Range("B5").Value = "No"
$ws.Range("B7").Value = "YES"
$ws.Range("B9").Value = "No"
$ws.Range("B10").Value = "No"
$ws.Range("B11").Value = "No"
$ws.Range("B14").Value = "No"
$ws.Range("B15").Value = "No"
$ws.Range("B16").Value = "No"
$ws.Range("B17").Value = "No"
$ws.Range("B18").Value = "No"
$ws.Range("B19").Value = "No"

# Update the view: scroll back to top-left A1 and move selection to B7
$ws.Range("A1").Select()
$ws.Range("B7").Select()
